$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row text (capitalisation / wording tweaks)
$ws.Range("A1").Value = "Nombre cliente"
$ws.Range("B1").Value = "Fecha de retiro"
$ws.Range("C1").Value = "Monto total"

# Give column B an explicit width (it had none before)
$ws.Columns.Item(2).ColumnWidth = 12.67

# Move the active selection like in the saved file
$ws.Range("E6").Select() | Out-Null
